$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "idetudiant" id column header ---
$ws.Range("A1").Value = "idetudiant"

# --- Header row (A1:H1) becomes bold ---
$ws.Range("A1:H1").Font.FontStyle = "Bold"

# --- Competence status highlighting ---
# Green fill (ADM / pass) cells
$ws.Range("G2:H2").Interior.Color = 65280
$ws.Range("G3").Interior.Color = 65280

# Red fill cells (everything else / empty competence cells)
$ws.Range("I2:K2").Interior.Color = 255
$ws.Range("H3:K3").Interior.Color = 255
$ws.Range("G4:K4").Interior.Color = 255
$ws.Range("G5:K5").Interior.Color = 255
$ws.Range("G6:K6").Interior.Color = 255

# --- Selection moved to K6 ---
$null = $ws.Range("K6").Select()
